# This script re-orders the per-observation data (id, taxon id/name/science
# name/author, coordinates, and observer names) across rows 3-17 of the
# active sheet. The mapping below says: the NEW data for row <key> is the
# OLD data that used to live in row <value> before this edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as a single "observation record".
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AW", "AX")

# Snapshot the current ("before") values for rows 3..17 so that writes to
# one row never clobber data we still need to read for another row.
$snapshot = @{}
for ($r = 3; $r -le 17; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Target row -> source row (source row's old values become the target row's
# new values).
$mapping = @{
    3  = 12
    4  = 3
    5  = 4
    6  = 5
    7  = 6
    8  = 13
    9  = 7
    10 = 14
    11 = 15
    12 = 8
    13 = 16
    14 = 9
    15 = 17
    16 = 10
    17 = 11
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value2 = $sourceVals[$col]
    }
}
